# LOQ4241.docx edit script
# Applies the changes described by the commit diff using Word COM-interop
# style calls against $word.ActiveDocument.

function Set-ItalicRunText($range, $text) {
    # Assign text to a (currently empty) Range, then italicize just the
    # text we inserted (exclude the trailing paragraph mark) so we don't
    # stamp rPr onto the paragraph mark itself.
    $range.Text = $text
    $s = $range.Start
    $e = $range.End
    $r2 = $d.Range($s, $e - 1)
    $r2.Italic = $true
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Work from the bottom of the document upward so paragraph insertions
# don't shift the indices of paragraphs we still need to touch.
# ---------------------------------------------------------------------

# --- Bibliografia: replace the whole (multi-run) paragraph with a
#     single new run of text, removing the old <w:br/> separated items.
$pBib = $d.Paragraphs(16)
$pBib.Range.Text = "FURTADO, N.; KAWAMOTO, E. Avaliação de Projetos de Transporte. São Carlos: Serviço Gráfico EESC-USP, 2002. 254 p.POWER, D. J. Decision Support Systems. London: Quorum Books, 2002. 251 p.GOMES, L. F. A. M.; GOMES, C. F. S.; ALMEIDA, A. T, Tomada de Decisão Gerencial: enfoque multicritério, São Paulo: Atlas, 2002.SHIMIZU, T., Decisão nas Organizações: introdução aos problemas de decisão encontrados nas organizações e nos sistemas de apoio à decisão, São Paulo: Atlas, 2001.DEVLIN, G. (ed.). Decision Support Systems: advances in. Zagreb: Intech, 2010. 342 p.GARCÍA-DÍAZ, V. Algorithms in Decision Support Systems. Basel: MDPI, 2020. 147 p."

# --- Avaliação: three inline text substitutions (Método / Critério /
#     Norma de recuperação), leaving the bold labels untouched.
$d.Content.Find.Execute("Aulas expositivas teóricas, aulas práticas, aulas de exercícios.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Provas, trabalhos em grupo, exercícios individuais e seminários.", 2)
$d.Content.Find.Execute("A Nota Final do aluno será determinada segundo a seguinte equação: Nota Final = (Prova- Bimestral-1*0,4) + (Prova-Bimestral-2*0,4) + (Trabalho*0,2)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Média das atividades avaliativas.", 2)
$d.Content.Find.Execute("Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação.", 2)

# --- Programa: replace the Portuguese body text, then add a new
#     italic English paragraph right after it.
$pPrograma = $d.Paragraphs(12)
$pPrograma.Range.Text = "i) Teoria da Decisão; ii) Estruturação de modelos de decisão; iii) Decisão com múltiplos cenários ou múltiplos critérios; iv) Decisão com incerteza; e v) Sistemas de auxílio à decisão e sistemas especialistas."
$pPrograma.Range.InsertParagraphAfter()
$pProgramaEn = $d.Paragraphs(13)
Set-ItalicRunText $pProgramaEn.Range "i) Decision Theory; ii) Decision support systems planning; iii) Decision with Multiple Scenarios or Multiple Criteria; iv) Decision with Uncertainty; and v) Decision Support Systems."

# --- Programa resumido: replace the Portuguese body text, then add a
#     new italic English paragraph right after it.
$pResumido = $d.Paragraphs(10)
$pResumido.Range.Text = "Teoria da Decisão; Planejamento de sistemas de apoio à decisão; Decisão com Múltiplos Cenários ou Múltiplos Critérios; Decisão com Incerteza; Sistemas de Auxílio à Decisão."
$pResumido.Range.InsertParagraphAfter()
$pResumidoEn = $d.Paragraphs(11)
Set-ItalicRunText $pResumidoEn.Range "Decision Theory; Decision support systems planning; Decision with Multiple Scenarios or Multiple Criteria; Decision with Uncertainty; Decision Support Systems."

# --- Docente(s) Responsável(eis): swap the listed professor.
$d.Content.Find.Execute("5840917 - Fabricio Maciel Gomes", $true, $false, $false, $false, $false,
                         $true, 1, $false, "3295113 - José Eduardo Holler Branco", 2)

# --- Objetivos: add a new italic English paragraph right after the
#     existing Portuguese objective paragraph.
$pObjetivos = $d.Paragraphs(6)
$pObjetivos.Range.InsertParagraphAfter()
$pObjetivosEn = $d.Paragraphs(7)
Set-ItalicRunText $pObjetivosEn.Range "Provide theory, tools and methods for supporting decision-making."

# --- Ativação date bump.
$d.Content.Find.Execute("Ativação: 01/01/2016", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ativação: 01/01/2024", 2)

# --- Subtitle (Heading3, previously empty): add the English title run.
$pSubtitle = $d.Paragraphs(2)
$pSubtitle.Range.Text = "Decision  Support Systems"

Write-Output "edit complete"
